# Agregada tabla independiente para referentes
#
# 1) Rename the existing sheet "Sheet1" -> "obras"
# 2) Add a new sheet "refentes" right after it, with its own
#    Título/Fecha/Periódico/Archivo table pulled out of "obras"
# 3) Fix a couple of stray values left in "obras" (typo'd filename,
#    and a referente title that had been mixed into the wrong cell)
# 4) Bold the header row on both sheets, zoom both to 150%, and leave
#    the selection / active sheet the way the author left them (on
#    "refentes").

$wb = $excel.ActiveWorkbook

# --- sheet 1: rename Sheet1 -> obras ------------------------------------
$obras = $wb.Worksheets.Item(1)
$obras.Name = "obras"

# --- fix the two stray cells in obras -----------------------------------
$obras.Range("I2").Value = "doble-suicidio-el-tiempo.jpg"
$obras.Range("I3").Value = "Una indígena y su hijo murieron en persecución"

# --- bold header row in obras --------------------------------------------
$obras.Range("A1:I1").Font.Bold = $true

# --- view settings for obras ----------------------------------------------
$obras.Application.ActiveWindow.Zoom = 150
$obras.Range("F3:I3").Select() | Out-Null

# --- sheet 2: add "refentes" right after "obras" --------------------------
$refentes = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $obras)
$refentes.Name = "refentes"

# --- column widths (characters) matching the authored layout -------------
$refentes.Columns.Item(1).ColumnWidth = 52.4987
$refentes.Columns.Item(2).ColumnWidth = 23.8307
$refentes.Columns.Item(3).ColumnWidth = 39.3307
$refentes.Columns.Item(4).ColumnWidth = 55.9987

# --- header row ------------------------------------------------------------
$refentes.Range("A1").Value = "Título"
$refentes.Range("B1").Value = "Fecha"
$refentes.Range("C1").Value = "Periódico"
$refentes.Range("D1").Value = "Archivo"
$refentes.Range("A1:D1").Font.Bold = $true

# --- data rows ---------------------------------------------------------------
$refentes.Range("A2").Value = 'Doble suicidio en "El Sisga"'
$refentes.Range("B2").Value = "Junio 29 1965"
$refentes.Range("C2").Value = "El Tiempo"
$refentes.Range("D2").Value = "doble-suicidio-el-tiempo.jpg"

$refentes.Range("A3").Value = "Una indígena y su hijo murieron en persecución"
$refentes.Range("B3").Value = "Mayo 24 del 96"
$refentes.Range("C3").Value = "El Tiempo"
$refentes.Range("D3").Value = "el-paraiso.jpg"

$refentes.Range("A4").Value = "Láminas de paisajes latinoamericanos"
$refentes.Range("D4").Value = "laminas-paisajes.jpg"

# --- view settings for refentes (this is the sheet left active/selected) --
$refentes.Application.ActiveWindow.Zoom = 150
$refentes.Range("C4").Select() | Out-Null
$refentes.Activate() | Out-Null
